# Auto-generated edit script applying scheduled-runner market data updates
# across worksheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# --- Worksheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H74").Value = 17907.072
$ws.Range("I74").Value = 19666.666
$ws.Range("J74").Value = 7349.5
$ws.Range("K74").Value = 19666.666
$ws.Range("L74").Value = 7349.5
$ws.Range("M74").Value = -18730.666
$ws.Range("N74").Value = -9221.5
$ws.Range("H77").Value = 17907.072
$ws.Range("I77").Value = 19666.666
$ws.Range("J77").Value = 7349.5
$ws.Range("K77").Value = 98333.33
$ws.Range("L77").Value = 36747.5
$ws.Range("M77").Value = -93653.33
$ws.Range("N77").Value = -46107.5
$ws.Range("H94").Value = 2355
$ws.Range("I94").Value = 1740.1666
$ws.Range("J94").Value = 4199.5
$ws.Range("K94").Value = 1740.1666
$ws.Range("L94").Value = 4199.5
$ws.Range("M94").Value = -1289.1666
$ws.Range("N94").Value = -5101.5
$ws.Range("H137").Value = 1473804.2
$ws.Range("I137").Value = 1788285
$ws.Range("J137").Value = 6227.3335
$ws.Range("K137").Value = 5364855
$ws.Range("L137").Value = 18682.0005
$ws.Range("M137").Value = -5362305
$ws.Range("N137").Value = -23782.0005
$ws.Range("H138").Value = 3397.279
$ws.Range("J138").Value = 3011.4412
$ws.Range("L138").Value = 9034.3236
$ws.Range("N138").Value = -19314.3236
$ws.Range("H141").Value = 2124.0833
$ws.Range("I141").Value = 1953.6364
$ws.Range("K141").Value = 5860.9092
$ws.Range("M141").Value = -680.9092000000001

# --- Worksheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 1896.6744
$ws.Range("I32").Value = 1906.8049
$ws.Range("K32").Value = 1906.8049
$ws.Range("M32").Value = -1619.8049
$ws.Range("H45").Value = 27858.875
$ws.Range("I45").Value = 36492.832
$ws.Range("K45").Value = 36492.832
$ws.Range("M45").Value = -36115.832
$ws.Range("H63").Value = 4743.1665
$ws.Range("I63").Value = 2364.75
$ws.Range("K63").Value = 2364.75
$ws.Range("M63").Value = -1678.75
$ws.Range("H66").Value = 4743.1665
$ws.Range("I66").Value = 2364.75
$ws.Range("K66").Value = 11823.75
$ws.Range("M66").Value = -8391.75
$ws.Range("H88").Value = 5761.857
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 5761.857
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 5761.857
$ws.Range("N88").Value = -6573.857
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 5761.857
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 5761.857
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 5761.857
$ws.Range("N91").Value = -8569.857
$ws.Range("M91").ClearContents()
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680
$ws.Range("H110").Value = 526
$ws.Range("I110").Value = 447.41666
$ws.Range("J110").Value = 997.5
$ws.Range("K110").Value = 447.41666
$ws.Range("L110").Value = 997.5
$ws.Range("M110").Value = 1597.58334
$ws.Range("N110").Value = -5087.5
$ws.Range("H122").Value = 3392.3333
$ws.Range("I122").Value = 3090.375
$ws.Range("J122").Value = 3633.9
$ws.Range("K122").Value = 9271.125
$ws.Range("L122").Value = 10901.7
$ws.Range("M122").Value = -6821.125
$ws.Range("N122").Value = -15801.7
$ws.Range("H129").Value = 25000
$ws.Range("J129").Value = 25000
$ws.Range("L129").Value = 25000
$ws.Range("N129").Value = -35000

# --- Worksheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 50008280
$ws.Range("I20").Value = 83344740
$ws.Range("K20").Value = 83344740
$ws.Range("M20").Value = -83344493
$ws.Range("H64").Value = 918.3333
$ws.Range("J64").Value = 1286.5
$ws.Range("L64").Value = 1286.5
$ws.Range("N64").Value = -1736.5
$ws.Range("H67").Value = 918.3333
$ws.Range("J67").Value = 1286.5
$ws.Range("L67").Value = 1286.5
$ws.Range("N67").Value = -2846.5
$ws.Range("H94").Value = 80001460
$ws.Range("J94").Value = 2954.1667
$ws.Range("L94").Value = 2954.1667
$ws.Range("N94").Value = -3856.1667
$ws.Range("H105").Value = 26002084
$ws.Range("I105").Value = 2001349.6
$ws.Range("J105").Value = 50002820
$ws.Range("K105").Value = 2001349.6
$ws.Range("L105").Value = 50002820
$ws.Range("M105").Value = -1999602.6
$ws.Range("N105").Value = -50006314
$ws.Range("H134").Value = 4755.091
$ws.Range("I134").Value = 5479.6
$ws.Range("K134").Value = 16438.8
$ws.Range("M134").Value = -13903.8

# --- Worksheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 1993.2222
$ws.Range("I22").Value = 1588.6
$ws.Range("J22").Value = 2499
$ws.Range("K22").Value = 1588.6
$ws.Range("L22").Value = 2499
$ws.Range("M22").Value = -1238.6
$ws.Range("N22").Value = -3199
$ws.Range("H107").Value = 2273550.8
$ws.Range("I107").Value = 3572048
$ws.Range("J107").Value = 1180.5
$ws.Range("K107").Value = 3572048
$ws.Range("L107").Value = 1180.5
$ws.Range("M107").Value = -3570128
$ws.Range("N107").Value = -5020.5

# --- Worksheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 1161
$ws.Range("I2").Value = 237.63637
$ws.Range("K2").Value = 1425.81822
$ws.Range("M2").Value = -1312.81822
$ws.Range("H4").Value = 11176267
$ws.Range("I4").Value = 11528052
$ws.Range("K4").Value = 34584156
$ws.Range("M4").Value = -34584044
$ws.Range("H5").Value = 605.6539
$ws.Range("J5").Value = 1449.5
$ws.Range("L5").Value = 4348.5
$ws.Range("N5").Value = -4572.5
$ws.Range("H47").Value = 2150.9
$ws.Range("I47").Value = 502.25
$ws.Range("K47").Value = 1506.75
$ws.Range("M47").Value = -1075.75
$ws.Range("H75").Value = 787.5
$ws.Range("J75").Value = 787.5
$ws.Range("L75").Value = 2362.5
$ws.Range("N75").Value = -4358.5
$ws.Range("H78").Value = 787.5
$ws.Range("J78").Value = 787.5
$ws.Range("L78").Value = 7087.5
$ws.Range("N78").Value = -17071.5
$ws.Range("H80").Value = 19727.727
$ws.Range("J80").Value = 23778.334
$ws.Range("L80").Value = 71335.00199999999
$ws.Range("N80").Value = -73207.00199999999
$ws.Range("H83").Value = 19727.727
$ws.Range("J83").Value = 23778.334
$ws.Range("L83").Value = 214005.006
$ws.Range("N83").Value = -223365.006
$ws.Range("H114").Value = 1518.5
$ws.Range("I114").Value = 1422.4
$ws.Range("K114").Value = 4267.200000000001
$ws.Range("M114").Value = -1013.200000000001
$ws.Range("H135").Value = 605.6539
$ws.Range("J135").Value = 1449.5
$ws.Range("L135").Value = 13045.5
$ws.Range("N135").Value = -18115.5
$ws.Range("H137").Value = 25003016
$ws.Range("I137").Value = 71431464
$ws.Range("J137").Value = 3082.3076
$ws.Range("K137").Value = 214294392
$ws.Range("L137").Value = 9246.9228
$ws.Range("M137").Value = -214289292
$ws.Range("N137").Value = -19446.9228

# --- Worksheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 10771.944
$ws.Range("I70").Value = 4738.154
$ws.Range("K70").Value = 4738.154
$ws.Range("M70").Value = -4468.154
$ws.Range("H73").Value = 10771.944
$ws.Range("I73").Value = 4738.154
$ws.Range("K73").Value = 4738.154
$ws.Range("M73").Value = -3802.154
$ws.Range("H80").Value = 333337150
$ws.Range("I80").Value = 500002500
$ws.Range("J80").Value = 6500
$ws.Range("K80").Value = 500002500
$ws.Range("L80").Value = 6500
$ws.Range("M80").Value = -500001502
$ws.Range("N80").Value = -8496
$ws.Range("H83").Value = 333337150
$ws.Range("I83").Value = 500002500
$ws.Range("J83").Value = 6500
$ws.Range("K83").Value = 2500012500
$ws.Range("L83").Value = 32500
$ws.Range("M83").Value = -2500007508
$ws.Range("N83").Value = -42484
$ws.Range("H102").Value = 2800.5
$ws.Range("I102").Value = 2281.9167
$ws.Range("K102").Value = 2281.9167
$ws.Range("M102").Value = -659.9167000000002

# --- Worksheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H82").Value = 1068.1818
$ws.Range("I82").Value = 1075.25
$ws.Range("J82").Value = 1049.3334
$ws.Range("K82").Value = 1075.25
$ws.Range("L82").Value = 1049.3334
$ws.Range("M82").Value = -714.25
$ws.Range("N82").Value = -1771.3334
$ws.Range("H85").Value = 1068.1818
$ws.Range("I85").Value = 1075.25
$ws.Range("J85").Value = 1049.3334
$ws.Range("K85").Value = 1075.25
$ws.Range("L85").Value = 1049.3334
$ws.Range("M85").Value = 172.75
$ws.Range("N85").Value = -3545.3334
$ws.Range("H100").Value = 6104
$ws.Range("I100").Value = 3901.2
$ws.Range("J100").Value = 7677.4287
$ws.Range("K100").Value = 3901.2
$ws.Range("L100").Value = 7677.4287
$ws.Range("M100").Value = -3360.2
$ws.Range("N100").Value = -8759.4287

# --- Worksheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 13159169
$ws.Range("I122").Value = 1324.0588
$ws.Range("K122").Value = 3972.1764
$ws.Range("M122").Value = -1522.1764
$ws.Range("H132").Value = 4205.607
$ws.Range("I132").Value = 3875.3076
$ws.Range("K132").Value = 11625.9228
$ws.Range("M132").Value = -9095.9228
$ws.Range("H136").Value = 244656.72
$ws.Range("I136").Value = 3741.4
$ws.Range("K136").Value = 11224.2
$ws.Range("M136").Value = -8674.200000000001
